$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest crypto quotes.
# We temporarily force a Text number format before writing each value so Excel
# does not reinterpret numeric-looking strings (e.g. "213.00", "6.90") as numbers
# and strip significant trailing zeros / dot separators, then restore the cells
# original style so formatting is left untouched.
$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.641.85"
$cell.Style = $origStyle
$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.34%  "
$cell.Style = $origStyle
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.637.27"
$cell.Style = $origStyle
$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.91%  "
$cell.Style = $origStyle
$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.09%  "
$cell.Style = $origStyle
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "213.00"
$cell.Style = $origStyle
$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.20%  "
$cell.Style = $origStyle
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.503"
$cell.Style = $origStyle
$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.39%  "
$cell.Style = $origStyle
$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.10%  "
$cell.Style = $origStyle
$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.34%  "
$cell.Style = $origStyle
$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.60%  "
$cell.Style = $origStyle
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.11"
$cell.Style = $origStyle
$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.01%  "
$cell.Style = $origStyle
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0841"
$cell.Style = $origStyle
$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.09%  "
$cell.Style = $origStyle
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.865.53"
$cell.Style = $origStyle
$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.90%  "
$cell.Style = $origStyle
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.628.09"
$cell.Style = $origStyle
$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell.Style = $origStyle
$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.74%  "
$cell.Style = $origStyle
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.528"
$cell.Style = $origStyle
$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.88%  "
$cell.Style = $origStyle
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.637.44"
$cell.Style = $origStyle
$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.29%  "
$cell.Style = $origStyle
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "63.43"
$cell.Style = $origStyle
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0745"
$cell.Style = $origStyle
$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.17%  "
$cell.Style = $origStyle
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "217.98"
$cell.Style = $origStyle
$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +7.46%  "
$cell.Style = $origStyle
$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.07%  "
$cell.Style = $origStyle
$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.48%  "
$cell.Style = $origStyle
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.52"
$cell.Style = $origStyle
$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.74%  "
$cell.Style = $origStyle
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.24"
$cell.Style = $origStyle
$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.35%  "
$cell.Style = $origStyle
$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.55%  "
$cell.Style = $origStyle
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "149.22"
$cell.Style = $origStyle
$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.73%  "
$cell.Style = $origStyle
$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.09%  "
$cell.Style = $origStyle
$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.11%  "
$cell.Style = $origStyle
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.90"
$cell.Style = $origStyle
$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.19%  "
$cell.Style = $origStyle
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.50"
$cell.Style = $origStyle
$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.70%  "
$cell.Style = $origStyle
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0522"
$cell.Style = $origStyle
$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.86%  "
$cell.Style = $origStyle
$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.49%  "
$cell.Style = $origStyle
$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.41%  "
$cell.Style = $origStyle
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.94"
$cell.Style = $origStyle
$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.30%  "
$cell.Style = $origStyle
$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.49%  "
$cell.Style = $origStyle
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.38"
$cell.Style = $origStyle
$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.58%  "
$cell.Style = $origStyle
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.182.16"
$cell.Style = $origStyle
$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.85%  "
$cell.Style = $origStyle
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0173"
$cell.Style = $origStyle
$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.87%  "
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.811"
$cell.Style = $origStyle
$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.07%  "
$cell.Style = $origStyle
$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.35%  "
$cell.Style = $origStyle
$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell.Style = $origStyle
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.31"
$cell.Style = $origStyle
$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.40%  "
$cell.Style = $origStyle
$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.78%  "
$cell.Style = $origStyle
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.793"
$cell.Style = $origStyle
$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.02%  "
$cell.Style = $origStyle
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.773.74"
$cell.Style = $origStyle
$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.82%  "
$cell.Style = $origStyle
$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.66%  "
$cell.Style = $origStyle
$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.48%  "
$cell.Style = $origStyle
$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.50%  "
$cell.Style = $origStyle
$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.97%  "
$cell.Style = $origStyle
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.64"
$cell.Style = $origStyle
$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.71%  "
$cell.Style = $origStyle
$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.19%  "
$cell.Style = $origStyle
$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.12%  "
$cell.Style = $origStyle
